$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Set E column ("num_commenti") to 100 for the rows that changed,
#        before any row shifting happens (rows 3-28 keep the same row
#        numbers throughout this script). ---
$eFixups = @{
    3  = 100
    4  = 100
    5  = 100
    6  = 100
    7  = 100
    8  = 100
    12 = 100
    13 = 100
    14 = 100
    15 = 100
    16 = 100
    17 = 100
    18 = 100
    19 = 100
    20 = 100
    21 = 100
    22 = 100
    23 = 100
    24 = 100
    25 = 100
    26 = 100
    27 = 100
    28 = 100
}
foreach ($r in $eFixups.Keys) {
    $ws.Cells.Item($r, 5).Value = $eFixups[$r]
}

# --- 2) Insert a brand new row at position 29: a new
#        "L'incoronazione di Re Carlo" / La Repubblica / YouTube entry.
#        This pushes the existing rows 29-37 ("La morte di Silvio
#        Berlusconi" block) down to rows 30-38. ---
$ws.Rows.Item(29).Insert()

$ws.Cells.Item(29, 1).Value = ""
$ws.Cells.Item(29, 2).Value = "L'incoronazione di Re Carlo"
$ws.Cells.Item(29, 3).Value = "La Repubblica"
$ws.Cells.Item(29, 4).Value = "YouTube"
$ws.Cells.Item(29, 5).Value = 100

# --- 3) Fix up the num_commenti values of the (now shifted) "La morte
#        di Silvio Berlusconi" rows (30-38). ---
$shiftedEFixups = @{
    30 = 100
    31 = 100
    32 = 100
    33 = 99
    34 = 100
    35 = 100
    36 = 1
    37 = 100
    38 = 100
}
foreach ($r in $shiftedEFixups.Keys) {
    $ws.Cells.Item($r, 5).Value = $shiftedEFixups[$r]
}

# --- 4) Within that shifted block, a few rows also changed which
#        giornale/social they reference (typo + reordering of the
#        remaining social-network rows). ---
$ws.Cells.Item(36, 3).Value = "Il Corrirere Della Sera"
$ws.Cells.Item(36, 4).Value = "Facebook"

$ws.Cells.Item(37, 3).Value = "La Repubblica"
$ws.Cells.Item(37, 4).Value = "Facebook"

$ws.Cells.Item(38, 3).Value = "La Repubblica"
$ws.Cells.Item(38, 4).Value = "Instagram"

# --- 5) Append a brand new final row 39: another "La morte di Silvio
#        Berlusconi" / La Repubblica / YouTube entry. ---
$ws.Cells.Item(39, 1).Value = ""
$ws.Cells.Item(39, 2).Value = "La morte di Silvio Berlusconi"

$ws.Cells.Item(39, 3).Value = "La Repubblica"
$ws.Cells.Item(39, 4).Value = "YouTube"
$ws.Cells.Item(39, 5).Value = 100
